$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 150 — 2024-07-30, an SMN.MI daily OHLC record with zero volume
# ---------------------------------------------------------------------------

# Column A needs the same date/time style as the rows above it (xf with the
# custom "yyyy-mm-dd hh:mm:ss" number format), so copy the style down from
# the previous row rather than building a new xf.
$ws.Range("A149").Copy()
$ws.Range("A150").PasteSpecial(-4122)
$ws.Range("A150").Value = 45503.2916666667

$ws.Range("B150").Value = 0
$ws.Range("C150").Value = 1.92999994754791
$ws.Range("D150").Value = 1.92999994754791
$ws.Range("E150").Value = 1.92999994754791
$ws.Range("F150").Value = 1.92999994754791

# Column G stores the adj_close as literal text (shared string), not a
# number. Driving it through TEXT()+paste-values keeps the cell a plain
# string without forcing a new "@" number-format style onto it.
$ws.Range("G150").Formula = "=TEXT(1.92999994754791,""0.00000000000000"")"
$ws.Range("G150").Copy()
$ws.Range("G150").PasteSpecial(-4163)

$ws.Range("H150").Value = "SMN.MI"

# ---------------------------------------------------------------------------
# New row 151 — 2024-07-31, an SMN.MI daily OHLC record with 3000 volume
# ---------------------------------------------------------------------------

$ws.Range("A149").Copy()
$ws.Range("A151").PasteSpecial(-4122)
$ws.Range("A151").Value = 45504.5192013889

$ws.Range("B151").Value = 3000
$ws.Range("C151").Value = 1.94000005722046
$ws.Range("D151").Value = 1.92999994754791
$ws.Range("E151").Value = 1.92999994754791
$ws.Range("F151").Value = 1.94000005722046

$ws.Range("G151").Formula = "=TEXT(1.94000005722046,""0.00000000000000"")"
$ws.Range("G151").Copy()
$ws.Range("G151").PasteSpecial(-4163)

$ws.Range("H151").Value = "SMN.MI"
